$wb = $excel.ActiveWorkbook
Write-Output $wb.ActiveSheet.Name
